$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 9262559
$ws.Range("I76").Value = 9262559
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 9262559
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -9262244
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 9262559
$ws.Range("I79").Value = 9262559
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 9262559
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -9261467
$ws.Range("N79").ClearContents()

$ws.Range("H111").Value = 1262.4286
$ws.Range("I111").Value = 984.25
$ws.Range("J111").Value = 1633.3334
$ws.Range("K111").Value = 2952.75
$ws.Range("L111").Value = 4900.0002
$ws.Range("M111").Value = 114.25
$ws.Range("N111").Value = -11034.0002

$ws.Range("H116").Value = 10173.77
$ws.Range("I116").Value = 10859.167
$ws.Range("K116").Value = 10859.167
$ws.Range("M116").Value = -7417.166999999999

$ws.Range("H132").Value = 1819.2142
$ws.Range("I132").Value = 1772.68
$ws.Range("J132").Value = 2207
$ws.Range("K132").Value = 5318.04
$ws.Range("L132").Value = 6621
$ws.Range("M132").Value = -2788.04
$ws.Range("N132").Value = -11681

$ws.Range("H137").Value = 1354.125
$ws.Range("J137").Value = 2309.9
$ws.Range("L137").Value = 6929.700000000001
$ws.Range("N137").Value = -12029.7

$ws.Range("H138").Value = 3601.3823
$ws.Range("I138").Value = 2481.276
$ws.Range("J138").Value = 4434.282
$ws.Range("K138").Value = 7443.828
$ws.Range("L138").Value = 13302.846
$ws.Range("M138").Value = -2303.828
$ws.Range("N138").Value = -23582.846

$ws.Range("H140").Value = 97200
$ws.Range("J140").Value = 97200
$ws.Range("L140").Value = 97200
$ws.Range("N140").Value = -107560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2885.524
$ws.Range("I61").Value = 2824.1
$ws.Range("K61").Value = 2824.1
$ws.Range("M61").Value = -2612.1

$ws.Range("H74").Value = 1330.375
$ws.Range("I74").Value = 1375.5714
$ws.Range("K74").Value = 1375.5714
$ws.Range("M74").Value = -501.5714

$ws.Range("H77").Value = 1330.375
$ws.Range("I77").Value = 1375.5714
$ws.Range("K77").Value = 6877.857
$ws.Range("M77").Value = -2509.857

$ws.Range("H132").Value = 2437.78
$ws.Range("I132").Value = 2104.6765
$ws.Range("J132").Value = 3145.625
$ws.Range("K132").Value = 6314.029500000001
$ws.Range("L132").Value = 9436.875
$ws.Range("M132").Value = -3784.029500000001
$ws.Range("N132").Value = -14496.875

$ws.Range("H136").Value = 2885.524
$ws.Range("I136").Value = 2824.1
$ws.Range("K136").Value = 8472.299999999999
$ws.Range("M136").Value = -5922.299999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 743.7646999999999
$ws.Range("I134").Value = 676.8570999999999
$ws.Range("J134").Value = 1056
$ws.Range("K134").Value = 2030.5713
$ws.Range("L134").Value = 3168
$ws.Range("M134").Value = 504.4287000000002
$ws.Range("N134").Value = -8238

$ws.Range("H138").Value = 50780
$ws.Range("J138").Value = 50780
$ws.Range("L138").Value = 50780
$ws.Range("N138").Value = -61060

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10126.855
$ws.Range("I31").Value = 3334.257
$ws.Range("J31").Value = 15925.415
$ws.Range("K31").Value = 3334.257
$ws.Range("L31").Value = 15925.415
$ws.Range("M31").Value = -3039.257
$ws.Range("N31").Value = -16515.415

$ws.Range("H34").Value = 10126.855
$ws.Range("I34").Value = 3334.257
$ws.Range("J34").Value = 15925.415
$ws.Range("K34").Value = 3334.257
$ws.Range("L34").Value = 15925.415
$ws.Range("M34").Value = -3132.257
$ws.Range("N34").Value = -16329.415

$ws.Range("H68").Value = 15935.167
$ws.Range("J68").Value = 15935.167
$ws.Range("L68").Value = 15935.167
$ws.Range("N68").Value = -17433.167

$ws.Range("H71").Value = 15935.167
$ws.Range("J71").Value = 15935.167
$ws.Range("L71").Value = 47805.501
$ws.Range("N71").Value = -55293.501

$ws.Range("H138").Value = 49960
$ws.Range("J138").Value = 49960
$ws.Range("L138").Value = 49960
$ws.Range("N138").Value = -60240

$ws.Range("H140").Value = 56667
$ws.Range("J140").Value = 56667
$ws.Range("L140").Value = 56667
$ws.Range("N140").Value = -67027

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 930.1414
$ws.Range("I113").Value = 617.1667
$ws.Range("J113").Value = 950.3333
$ws.Range("K113").Value = 1851.5001
$ws.Range("L113").Value = 2850.9999
$ws.Range("M113").Value = 318.4999
$ws.Range("N113").Value = -7190.9999

$ws.Range("H137").Value = 7081.381
$ws.Range("I137").Value = 7039.3887
$ws.Range("J137").Value = 7333.3335
$ws.Range("K137").Value = 21118.1661
$ws.Range("L137").Value = 22000.0005
$ws.Range("M137").Value = -16018.1661
$ws.Range("N137").Value = -32200.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 691073.3
$ws.Range("I11").Value = 2575025
$ws.Range("K11").Value = 2575025
$ws.Range("M11").Value = -2574886

$ws.Range("H21").Value = 2500850
$ws.Range("I21").Value = 5000200
$ws.Range("K21").Value = 5000200
$ws.Range("M21").Value = -5000027

$ws.Range("H24").Value = 1500
$ws.Range("J24").Value = 1500
$ws.Range("L24").Value = 1500
$ws.Range("N24").Value = -1846

$ws.Range("H30").Value = 2500850
$ws.Range("I30").Value = 5000200
$ws.Range("K30").Value = 5000200
$ws.Range("M30").Value = -5000095

$ws.Range("H126").Value = 33335114
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 33335114
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 100005342
$ws.Range("N126").Value = -100010282
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 1710.0244
$ws.Range("I132").Value = 1307.6333
$ws.Range("J132").Value = 2807.4546
$ws.Range("K132").Value = 3922.8999
$ws.Range("L132").Value = 8422.363799999999
$ws.Range("M132").Value = -1392.8999
$ws.Range("N132").Value = -13482.3638

$ws.Range("H139").Value = 54075.332
$ws.Range("J139").Value = 54075.332
$ws.Range("L139").Value = 54075.332
$ws.Range("N139").Value = -64355.332

$ws.Range("H140").Value = 75780
$ws.Range("J140").Value = 75780
$ws.Range("L140").Value = 75780
$ws.Range("N140").Value = -86140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4167.923
$ws.Range("I7").Value = 3547.1667
$ws.Range("J7").Value = 4700
$ws.Range("K7").Value = 3547.1667
$ws.Range("L7").Value = 4700
$ws.Range("M7").Value = -3435.1667
$ws.Range("N7").Value = -4924

$ws.Range("H82").Value = 2118.0625
$ws.Range("I82").Value = 2381.125
$ws.Range("J82").Value = 1855
$ws.Range("K82").Value = 2381.125
$ws.Range("L82").Value = 1855
$ws.Range("M82").Value = -2020.125
$ws.Range("N82").Value = -2577

$ws.Range("H85").Value = 2118.0625
$ws.Range("I85").Value = 2381.125
$ws.Range("J85").Value = 1855
$ws.Range("K85").Value = 2381.125
$ws.Range("L85").Value = 1855
$ws.Range("M85").Value = -1133.125
$ws.Range("N85").Value = -4351

$ws.Range("H126").Value = 4167.923
$ws.Range("I126").Value = 3547.1667
$ws.Range("J126").Value = 4700
$ws.Range("K126").Value = 10641.5001
$ws.Range("L126").Value = 14100
$ws.Range("M126").Value = -8171.500100000001
$ws.Range("N126").Value = -19040

$ws.Range("H136").Value = 4332.04
$ws.Range("I136").Value = 2521.5293
$ws.Range("J136").Value = 8179.375
$ws.Range("K136").Value = 7564.5879
$ws.Range("L136").Value = 24538.125
$ws.Range("M136").Value = -5014.5879
$ws.Range("N136").Value = -29638.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 28533.334
$ws.Range("J99").Value = 28533.334
$ws.Range("L99").Value = 28533.334
$ws.Range("N99").Value = -34523.334

$ws.Range("H126").Value = 2813.3333
$ws.Range("I126").Value = 2629.6667
$ws.Range("J126").Value = 2997
$ws.Range("K126").Value = 7889.000100000001
$ws.Range("L126").Value = 8991
$ws.Range("M126").Value = -5419.000100000001
$ws.Range("N126").Value = -13931

$ws.Range("H139").Value = 57325
$ws.Range("J139").Value = 57325
$ws.Range("L139").Value = 57325
$ws.Range("N139").Value = -67605
Write-Host "Applied Sheets data refresh updates."
